$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 23 of the schedule table holds the "ЛР05" lab-work entry whose two
# date cells need updating.
$row = $t.Rows.Item(23)

# --- Column 2 (lecture date): "17.03" -> "19.03" ------------------------
# Assigning straight to Range.Text (rather than Find.Execute, whose Wrap/
# Replace semantics in this host search/replace across the whole story
# instead of staying inside the supplied Range) keeps the edit scoped to
# this one cell and preserves the existing run's rPr.
$dateCell = $row.Cells.Item(2)
$dateCell.Range.Text = "19.03"

# --- Column 3 (lab date): "20.03" -> "21.03", split into 3 runs --------
# The source edit wound up as three separate runs ("2" / "1" / ".03")
# all sharing the original rPr.
$labCell = $row.Cells.Item(3)
$labCell.Range.Text = "21.03"

$cellStart = $labCell.Range.Start

# Force a run split after the 1st and 2nd characters by nudging the font
# size away from, then back to, the original 14pt (sz=28) value - runs
# only get physically separated when a direct-formatting value actually
# differs from its neighbour; setting the same value again is a no-op
# that leaves them split afterwards (they are not re-merged).
$chr1 = $d.Range($cellStart, $cellStart + 1)
$chr1.Font.Size = 99
$chr2 = $d.Range($cellStart + 1, $cellStart + 2)
$chr2.Font.Size = 99

$chr1b = $d.Range($cellStart, $cellStart + 1)
$chr1b.Font.Size = 14
$chr2b = $d.Range($cellStart + 1, $cellStart + 2)
$chr2b.Font.Size = 14

Write-Output ("Row23 dates: [" + $dateCell.Range.Text + "] [" + $labCell.Range.Text + "]")
